# Update the HLS stream links for Willow TV and Sky sports with refreshed
# md5/expires query parameters, then leave the selection on the Sky sports row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "https://off1.dunyapurkaraja.com:1686/hls/willowusa.m3u8?md5=DrujHrZ7SR-kyVZ9dTBnhg&expires=1742192590"
$ws.Range("B3").Value = "https://off1.dunyapurkaraja.com:1686/hls/skyscric.m3u8?md5=jgLvWH4aBWP6qAcnDQAclw&expires=1742192627"

$ws.Range("B3").Select()
